$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 241
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 4859
$ws.Range("F5").Value = 212
$ws.Range("F6").Value = 165
$ws.Range("F7").Value = 126
$ws.Range("F8").Value = 116
$ws.Range("F10").Value = 769
$ws.Range("F11").Value = 238
$ws.Range("F12").Value = 1220
$ws.Range("F13").Value = 124
$ws.Range("F16").Value = 89
$ws.Range("F20").Value = 4133
$ws.Range("F21").Value = 6440
$ws.Range("F24").Value = 88
$ws.Range("F27").Value = 4011
$ws.Range("F28").Value = 413
$ws.Range("F29").Value = 54
$ws.Range("F30").Value = 31
$ws.Range("F31").Value = 2611
$ws.Range("F34").Value = 153
$ws.Range("F35").Value = 311
$ws.Range("F37").Value = 382
$ws.Range("F38").Value = 190
$ws.Range("F40").Value = 1579
$ws.Range("F43").Value = 82
$ws.Range("F44").Value = 60
$ws.Range("F45").Value = 504
$ws.Range("F47").Value = 6
$ws.Range("F48").Value = 80

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 112

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 241
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 4859
$ws.Range("F6").Value = 165
$ws.Range("F9").Value = 116
$ws.Range("F12").Value = 238
$ws.Range("F13").Value = 1220
$ws.Range("F14").Value = 124
$ws.Range("F15").Value = 201
$ws.Range("F16").Value = 89
$ws.Range("F20").Value = 4133
$ws.Range("F21").Value = 6440
$ws.Range("F22").Value = 39
$ws.Range("F24").Value = 88
$ws.Range("F25").Value = 545
$ws.Range("F26").Value = 48
$ws.Range("F27").Value = 4011
$ws.Range("F28").Value = 413
$ws.Range("F29").Value = 54
$ws.Range("F30").Value = 31
$ws.Range("F31").Value = 2611
$ws.Range("F35").Value = 311
$ws.Range("F37").Value = 382
$ws.Range("F38").Value = 190
$ws.Range("F40").Value = 1579
$ws.Range("F41").Value = 980
$ws.Range("F42").Value = 49
$ws.Range("F43").Value = 82
$ws.Range("F44").Value = 60
$ws.Range("F46").Value = 485
$ws.Range("F49").Value = 596
